# This script updates the "Pais" COVID-19 country statistics sheet to a newer
# data snapshot (refresh timestamp 19:28 instead of 18:11). The underlying feed
# refreshed totals for many countries; for a handful of countries whose running
# totals crossed one another in the new snapshot, the row order changes too
# (e.g. Israel now outranks Ucrania, Libano now outranks Bosnia y Herzegovina,
# Maldivas now outranks Tayikistan, Mozambique moves up, Trinidad y Tobago moves
# up, and Reunion now outranks Estonia). Only the rows that actually changed
# value and/or country name are touched; everything else is left as-is.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value2 = "Datos actualizados a 12 de Septiembre de 2020 a las 19:28"

# Row 4: Estados Unidos
$ws.Range("A4").Value2 = "Estados Unidos"
$ws.Range("B4").Value2 = 6652721
$ws.Range("C4").Value2 = 16474
$ws.Range("D4").Value2 = 3927057
$ws.Range("E4").Value2 = 2527986
$ws.Range("F4").Value2 = 0
$ws.Range("G4").Value2 = 257
$ws.Range("H4").Value2 = 197678

# Row 5: India
$ws.Range("A5").Value2 = "India"
$ws.Range("B5").Value2 = 4742743
$ws.Range("C5").Value2 = 85364
$ws.Range("D5").Value2 = 3693206
$ws.Range("E5").Value2 = 970985
$ws.Range("F5").Value2 = 0
$ws.Range("G5").Value2 = 1046
$ws.Range("H5").Value2 = 78552

# Row 6: Brasil
$ws.Range("A6").Value2 = "Brasil"
$ws.Range("B6").Value2 = 4297949
$ws.Range("C6").Value2 = 13971
$ws.Range("D6").Value2 = 3530655
$ws.Range("E6").Value2 = 636424
$ws.Range("F6").Value2 = 0
$ws.Range("G6").Value2 = 396
$ws.Range("H6").Value2 = 130870

# Row 21: Turquia
$ws.Range("A21").Value2 = "Turquia"
$ws.Range("B21").Value2 = 289635
$ws.Range("C21").Value2 = 1509
$ws.Range("D21").Value2 = 257731
$ws.Range("E21").Value2 = 24905
$ws.Range("F21").Value2 = 0
$ws.Range("G21").Value2 = 48
$ws.Range("H21").Value2 = 6999

# Row 24: Alemania
$ws.Range("A24").Value2 = "Alemania"
$ws.Range("B24").Value2 = 260277
$ws.Range("C24").Value2 = 552
$ws.Range("D24").Value2 = 234850
$ws.Range("E24").Value2 = 16003
$ws.Range("F24").Value2 = 0
$ws.Range("G24").Value2 = 1
$ws.Range("H24").Value2 = 9424

# Row 27: Israel
$ws.Range("A27").Value2 = "Israel"
$ws.Range("B27").Value2 = 152525
$ws.Range("C27").Value2 = 3961
$ws.Range("D27").Value2 = 113494
$ws.Range("E27").Value2 = 37930
$ws.Range("F27").Value2 = 0
$ws.Range("G27").Value2 = 11
$ws.Range("H27").Value2 = 1101

# Row 28: Ucrania
$ws.Range("A28").Value2 = "Ucrania"
$ws.Range("B28").Value2 = 151859
$ws.Range("C28").Value2 = 3103
$ws.Range("D28").Value2 = 68346
$ws.Range("E28").Value2 = 80365
$ws.Range("F28").Value2 = 0
$ws.Range("G28").Value2 = 72
$ws.Range("H28").Value2 = 3148

# Row 29: Canada
$ws.Range("A29").Value2 = "Canada"
$ws.Range("B29").Value2 = 136102
$ws.Range("C29").Value2 = 476
$ws.Range("D29").Value2 = 120004
$ws.Range("E29").Value2 = 6928
$ws.Range("F29").Value2 = 0
$ws.Range("G29").Value2 = 7
$ws.Range("H29").Value2 = 9170

# Row 43: Marruecos
$ws.Range("A43").Value2 = "Marruecos"
$ws.Range("B43").Value2 = 84435
$ws.Range("C43").Value2 = 2238
$ws.Range("D43").Value2 = 65867
$ws.Range("E43").Value2 = 17015
$ws.Range("F43").Value2 = 0
$ws.Range("G43").Value2 = 29
$ws.Range("H43").Value2 = 1553

# Row 51: Etiopia
$ws.Range("A51").Value2 = "Etiopia"
$ws.Range("B51").Value2 = 63888
$ws.Range("C51").Value2 = 521
$ws.Range("D51").Value2 = 24493
$ws.Range("E51").Value2 = 38399
$ws.Range("F51").Value2 = 0
$ws.Range("G51").Value2 = 10
$ws.Range("H51").Value2 = 996

# Row 59: Argelia
$ws.Range("A59").Value2 = "Argelia"
$ws.Range("B59").Value2 = 48007
$ws.Range("C59").Value2 = 255
$ws.Range("D59").Value2 = 33875
$ws.Range("E59").Value2 = 12527
$ws.Range("F59").Value2 = 0
$ws.Range("G59").Value2 = 6
$ws.Range("H59").Value2 = 1605

# Row 63: Ghana
$ws.Range("A63").Value2 = "Ghana"
$ws.Range("B63").Value2 = 45434
$ws.Range("C63").Value2 = 46
$ws.Range("D63").Value2 = 44342
$ws.Range("E63").Value2 = 806
$ws.Range("F63").Value2 = 0
$ws.Range("G63").Value2 = 1
$ws.Range("H63").Value2 = 286

# Row 68: Kenia
$ws.Range("A68").Value2 = "Kenia"
$ws.Range("B68").Value2 = 35969
$ws.Range("C68").Value2 = 176
$ws.Range("D68").Value2 = 22771
$ws.Range("E68").Value2 = 12579
$ws.Range("F68").Value2 = 0
$ws.Range("G68").Value2 = 3
$ws.Range("H68").Value2 = 619

# Row 69: Chequia
$ws.Range("A69").Value2 = "Chequia"
$ws.Range("B69").Value2 = 34744
$ws.Range("C69").Value2 = 884
$ws.Range("D69").Value2 = 21150
$ws.Range("E69").Value2 = 13141
$ws.Range("F69").Value2 = 0
$ws.Range("G69").Value2 = 3
$ws.Range("H69").Value2 = 453

# Row 72: Irlanda
$ws.Range("A72").Value2 = "Irlanda"
$ws.Range("B72").Value2 = 30730
$ws.Range("C72").Value2 = 159
$ws.Range("D72").Value2 = 23364
$ws.Range("E72").Value2 = 5583
$ws.Range("F72").Value2 = 0
$ws.Range("G72").Value2 = 2
$ws.Range("H72").Value2 = 1783

# Row 77: Libano
$ws.Range("A77").Value2 = "Libano"
$ws.Range("B77").Value2 = 23669
$ws.Range("C77").Value2 = 686
$ws.Range("D77").Value2 = 7312
$ws.Range("E77").Value2 = 16118
$ws.Range("F77").Value2 = 0
$ws.Range("G77").Value2 = 10
$ws.Range("H77").Value2 = 239

# Row 78: Bosnia y Herzegovina
$ws.Range("A78").Value2 = "Bosnia y Herzegovina"
$ws.Range("B78").Value2 = 23138
$ws.Range("C78").Value2 = 304
$ws.Range("D78").Value2 = 15922
$ws.Range("E78").Value2 = 6526
$ws.Range("F78").Value2 = 0
$ws.Range("G78").Value2 = 4
$ws.Range("H78").Value2 = 690

# Row 91: Grecia
$ws.Range("A91").Value2 = "Grecia"
$ws.Range("B91").Value2 = 13036
$ws.Range("C91").Value2 = 302
$ws.Range("D91").Value2 = 3804
$ws.Range("E91").Value2 = 8930
$ws.Range("F91").Value2 = 0
$ws.Range("G91").Value2 = 2
$ws.Range("H91").Value2 = 302

# Row 94: Albania
$ws.Range("A94").Value2 = "Albania"
$ws.Range("B94").Value2 = 11185
$ws.Range("C94").Value2 = 164
$ws.Range("D94").Value2 = 6494
$ws.Range("E94").Value2 = 4361
$ws.Range("F94").Value2 = 0
$ws.Range("G94").Value2 = 3
$ws.Range("H94").Value2 = 330

# Row 100: Maldivas
$ws.Range("A100").Value2 = "Maldivas"
$ws.Range("B100").Value2 = 9052
$ws.Range("C100").Value2 = 62
$ws.Range("D100").Value2 = 7055
$ws.Range("E100").Value2 = 1966
$ws.Range("F100").Value2 = 0
$ws.Range("G100").Value2 = 0
$ws.Range("H100").Value2 = 31

# Row 101: Tayikistan
$ws.Range("A101").Value2 = "Tayikistan"
$ws.Range("B101").Value2 = 9014
$ws.Range("C101").Value2 = 37
$ws.Range("D101").Value2 = 7782
$ws.Range("E101").Value2 = 1160
$ws.Range("F101").Value2 = 0
$ws.Range("G101").Value2 = 0
$ws.Range("H101").Value2 = 72

# Row 113: Mozambique
$ws.Range("A113").Value2 = "Mozambique"
$ws.Range("B113").Value2 = 5040
$ws.Range("C113").Value2 = 122
$ws.Range("D113").Value2 = 2905
$ws.Range("E113").Value2 = 2100
$ws.Range("F113").Value2 = 0
$ws.Range("G113").Value2 = 4
$ws.Range("H113").Value2 = 35

# Row 114: Suazilandia
$ws.Range("A114").Value2 = "Suazilandia"
$ws.Range("B114").Value2 = 5025
$ws.Range("C114").Value2 = 0
$ws.Range("D114").Value2 = 4165
$ws.Range("E114").Value2 = 762
$ws.Range("F114").Value2 = 0
$ws.Range("G114").Value2 = 0
$ws.Range("H114").Value2 = 98

# Row 115: Guinea Ecuatorial
$ws.Range("A115").Value2 = "Guinea Ecuatorial"
$ws.Range("B115").Value2 = 4996
$ws.Range("C115").Value2 = 0
$ws.Range("D115").Value2 = 4490
$ws.Range("E115").Value2 = 423
$ws.Range("F115").Value2 = 0
$ws.Range("G115").Value2 = 0
$ws.Range("H115").Value2 = 83

# Row 116: Hong Kong
$ws.Range("A116").Value2 = "Hong Kong"
$ws.Range("B116").Value2 = 4939
$ws.Range("C116").Value2 = 13
$ws.Range("D116").Value2 = 4613
$ws.Range("E116").Value2 = 226
$ws.Range("F116").Value2 = 0
$ws.Range("G116").Value2 = 1
$ws.Range("H116").Value2 = 100

# Row 117: Congo
$ws.Range("A117").Value2 = "Congo"
$ws.Range("B117").Value2 = 4928
$ws.Range("C117").Value2 = 0
$ws.Range("D117").Value2 = 3887
$ws.Range("E117").Value2 = 953
$ws.Range("F117").Value2 = 0
$ws.Range("G117").Value2 = 0
$ws.Range("H117").Value2 = 88

# Row 134: Sri Lanka
$ws.Range("A134").Value2 = "Sri Lanka"
$ws.Range("B134").Value2 = 3195
$ws.Range("C134").Value2 = 26
$ws.Range("D134").Value2 = 2983
$ws.Range("E134").Value2 = 200
$ws.Range("F134").Value2 = 0
$ws.Range("G134").Value2 = 0
$ws.Range("H134").Value2 = 12

# Row 137: Mali
$ws.Range("A137").Value2 = "Mali"
$ws.Range("B137").Value2 = 2916
$ws.Range("C137").Value2 = 4
$ws.Range("D137").Value2 = 2276
$ws.Range("E137").Value2 = 512
$ws.Range("F137").Value2 = 0
$ws.Range("G137").Value2 = 0
$ws.Range("H137").Value2 = 128

# Row 139: Trinidad yTobago
$ws.Range("A139").Value2 = "Trinidad yTobago"
$ws.Range("B139").Value2 = 2892
$ws.Range("C139").Value2 = 67
$ws.Range("D139").Value2 = 766
$ws.Range("E139").Value2 = 2076
$ws.Range("F139").Value2 = 0
$ws.Range("G139").Value2 = 0
$ws.Range("H139").Value2 = 50

# Row 140: Bahamas
$ws.Range("A140").Value2 = "Bahamas"
$ws.Range("B140").Value2 = 2874
$ws.Range("C140").Value2 = 60
$ws.Range("D140").Value2 = 1285
$ws.Range("E140").Value2 = 1522
$ws.Range("F140").Value2 = 0
$ws.Range("G140").Value2 = 2
$ws.Range("H140").Value2 = 67

# Row 141: Reunion
$ws.Range("A141").Value2 = "Reunion"
$ws.Range("B141").Value2 = 2723
$ws.Range("C141").Value2 = 100
$ws.Range("D141").Value2 = 1313
$ws.Range("E141").Value2 = 1396
$ws.Range("F141").Value2 = 0
$ws.Range("G141").Value2 = 0
$ws.Range("H141").Value2 = 14

# Row 142: Estonia
$ws.Range("A142").Value2 = "Estonia"
$ws.Range("B142").Value2 = 2655
$ws.Range("C142").Value2 = 23
$ws.Range("D142").Value2 = 2252
$ws.Range("E142").Value2 = 339
$ws.Range("F142").Value2 = 0
$ws.Range("G142").Value2 = 0
$ws.Range("H142").Value2 = 64

# Row 144: Sudan del Sur
$ws.Range("A144").Value2 = "Sudan del Sur"
$ws.Range("B144").Value2 = 2578
$ws.Range("C144").Value2 = 10
$ws.Range("D144").Value2 = 1290
$ws.Range("E144").Value2 = 1239
$ws.Range("F144").Value2 = 0
$ws.Range("G144").Value2 = 0
$ws.Range("H144").Value2 = 49

# Row 150: Sierra Leona
$ws.Range("A150").Value2 = "Sierra Leona"
$ws.Range("B150").Value2 = 2096
$ws.Range("C150").Value2 = 9
$ws.Range("D150").Value2 = 1634
$ws.Range("E150").Value2 = 390
$ws.Range("F150").Value2 = 0
$ws.Range("G150").Value2 = 0
$ws.Range("H150").Value2 = 72

# Row 162: Liberia
$ws.Range("A162").Value2 = "Liberia"
$ws.Range("B162").Value2 = 1316
$ws.Range("C162").Value2 = 1
$ws.Range("D162").Value2 = 1210
$ws.Range("E162").Value2 = 24
$ws.Range("F162").Value2 = 0
$ws.Range("G162").Value2 = 0
$ws.Range("H162").Value2 = 82

# Row 179: Islas Feroe
$ws.Range("A179").Value2 = "Islas Feroe"
$ws.Range("B179").Value2 = 418
$ws.Range("C179").Value2 = 2
$ws.Range("D179").Value2 = 410
$ws.Range("E179").Value2 = 8
$ws.Range("F179").Value2 = 0
$ws.Range("G179").Value2 = 0
$ws.Range("H179").Value2 = 0

# Row 195: Liechtenstein
$ws.Range("A195").Value2 = "Liechtenstein"
$ws.Range("B195").Value2 = 111
$ws.Range("C195").Value2 = 2
$ws.Range("D195").Value2 = 105
$ws.Range("E195").Value2 = 5
$ws.Range("F195").Value2 = 0
$ws.Range("G195").Value2 = 0
$ws.Range("H195").Value2 = 1
